# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Chequia" / "Banglades" labels (row 45 <-> row 46) ---
# Before: A45 = Chequia, A46 = Banglades
# After:  A45 = Banglades, A46 = Chequia
$ws.Cells.Item(45, 1).Value = "Banglades"
$ws.Cells.Item(46, 1).Value = "Chequia"

# --- Update the "Datos actualizados" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 10:52"

# --- Row 34 (Polonia): update Casos activos / Recuperados ---
$ws.Cells.Item(34, 4).Value = 3236
$ws.Cells.Item(34, 5).Value = 8917

# --- Row 36 (Rumania): update Recuperados / Muertes hoy / Muertes ---
$ws.Cells.Item(36, 5).Value = 7714
$ws.Cells.Item(36, 7).Value = 2
$ws.Cells.Item(36, 8).Value = 695

# --- Row 45 (now Banglades): updated stats ---
$ws.Cells.Item(45, 2).Value = 7667
$ws.Cells.Item(45, 3).Value = 564
$ws.Cells.Item(45, 4).Value = 160
$ws.Cells.Item(45, 5).Value = 7339
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(45, 7).Value = 5
$ws.Cells.Item(45, 8).Value = 168

# --- Row 46 (now Chequia): updated stats (previously Chequia's old values) ---
$ws.Cells.Item(46, 2).Value = 7581
$ws.Cells.Item(46, 3).Value = 2
$ws.Cells.Item(46, 4).Value = 3120
$ws.Cells.Item(46, 5).Value = 4234
$ws.Cells.Item(46, 6).Value = 68
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 227

# --- Row 51 (Malasia): update Casos totales / Nuevos casos / Casos activos / Recuperados / Muertes hoy / Muertes ---
$ws.Cells.Item(51, 2).Value = 6002
$ws.Cells.Item(51, 3).Value = 57
$ws.Cells.Item(51, 4).Value = 4171
$ws.Cells.Item(51, 5).Value = 1729
$ws.Cells.Item(51, 7).Value = 2
$ws.Cells.Item(51, 8).Value = 102
